$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting for numeric-looking values
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.003.91'
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("D3").Value = '1.909.04'
$ws.Range("E3").Value = '  -3.99%  '
$ws.Range("D4").Value = '1.004'
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").Value = '325.00'
$ws.Range("E5").Value = '  -1.43%  '
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '0.4606'
$ws.Range("E7").Value = '  -1.72%  '
$ws.Range("D8").Value = '0.3834'
$ws.Range("E8").Value = '  -2.89%  '
$ws.Range("D9").Value = '0.07754'
$ws.Range("E9").Value = '  -2.79%  '
$ws.Range("D10").Value = '0.9863'
$ws.Range("E10").Value = '  -1.70%  '
$ws.Range("D11").Value = '22.11'
$ws.Range("E11").Value = '  -3.01%  '
$ws.Range("D12").Value = '1.927.54'
$ws.Range("E12").Value = '  -3.23%  '
$ws.Range("D13").Value = '7.006'
$ws.Range("E13").Value = '  -3.71%  '
$ws.Range("D14").Value = '5.710'
$ws.Range("E14").Value = '  -3.06%  '
$ws.Range("D15").Value = '0.07078'
$ws.Range("E15").Value = '  -1.30%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  -0.28%  '
$ws.Range("D17").Value = '84.11'
$ws.Range("E17").Value = '  -5.51%  '
$ws.Range("D18").Value = '0.000009575'
$ws.Range("E18").Value = '  -4.01%  '
$ws.Range("D19").Value = '16.79'
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("D20").Value = '1.002'
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("D21").Value = '29.031.80'
$ws.Range("E21").Value = '  -2.02%  '
$ws.Range("D22").Value = '5.341'
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("D23").Value = '10.98'
$ws.Range("E23").Value = '  -2.88%  '
$ws.Range("D24").Value = '2.146.33'
$ws.Range("E24").Value = '  -4.75%  '
$ws.Range("D25").Value = '2.079'
$ws.Range("E25").Value = '  -2.18%  '
$ws.Range("D26").Value = '156.59'
$ws.Range("E26").Value = '  -1.07%  '
$ws.Range("D27").Value = '19.24'
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("D28").Value = '5.620'
$ws.Range("E28").Value = '  -6.29%  '
$ws.Range("D29").Value = '118.13'
$ws.Range("D30").Value = '1.835'
$ws.Range("E30").Value = '  -6.67%  '
$ws.Range("D31").Value = '0.09285'
$ws.Range("E31").Value = '  -1.79%  '
$ws.Range("D32").Value = '0.8660'
$ws.Range("E32").Value = '  -3.08%  '
$ws.Range("D33").Value = '5.121'
$ws.Range("E33").Value = '  -3.27%  '
$ws.Range("D34").Value = '1.256'
$ws.Range("E34").Value = '  -6.86%  '
$ws.Range("D35").Value = '3.019'
$ws.Range("E35").Value = '  -5.43%  '
$ws.Range("D36").Value = '0.05749'
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("D37").Value = '1.153'
$ws.Range("E37").Value = '  -2.13%  '
$ws.Range("D38").Value = '1.003'
$ws.Range("E38").Value = '  -0.28%  '
$ws.Range("D39").Value = '0.02054'
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("D40").Value = '7.506'
$ws.Range("E40").Value = '  -5.32%  '
$ws.Range("D41").Value = '0.5548'
$ws.Range("E41").Value = '  -3.80%  '
$ws.Range("D42").Value = '0.1762'
$ws.Range("D43").Value = '9.316'
$ws.Range("E43").Value = '  -5.31%  '
$ws.Range("D44").Value = '2.737'
$ws.Range("E44").Value = '  +3.37%  '
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").Value = '0.5229'
$ws.Range("E45").Value = '  -2.98%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '11.30'
$ws.Range("E46").Value = '  -6.99%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '2.117'
$ws.Range("E47").Value = '  -2.26%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.06828'
$ws.Range("E48").Value = '  -2.02%  '
$ws.Range("B49").Value = 'PEPE'
$ws.Range("C49").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D49").Value = '0.000002611'
$ws.Range("E49").Value = '  -15.48%  '
$ws.Range("D50").Value = '112.15'
$ws.Range("D51").Value = '1.781'
$ws.Range("E51").Value = '  -4.89%  '
